$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.46"
$ws.Range("D3").Value = "'22.81"
$ws.Range("D4").Value = "'5.467"
$ws.Range("D5").Value = "'0.05755"
$ws.Range("D6").Value = "'3.431"
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = "'0.8121"
$ws.Range("E7").Value = '6MXTokenMX'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = "'0.8847"
$ws.Range("E8").Value = '7FTXTokenFTT'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'0.1442"
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = "'0.07329"
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = "'0.03011"
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.03132"
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09402"
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Value = "'3.930"
$ws.Range("E14").Value = '13MCDexMCB'
$ws.Range("D15").Value = "'0.001584"
$ws.Range("D16").Value = "'0.04813"
$ws.Range("D17").Value = "'0.0005838"
$ws.Range("D18").Value = "'0.006159"
$ws.Range("D19").Value = "'0.005125"
$ws.Range("D20").Value = "'0.0009971"
$ws.Range("D21").Value = "'0.0001499"
$ws.Range("B23").Value = 'KuCoinToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D23").Value = "'6.331"
$ws.Range("E23").Value = '22KuCoinTokenKCS'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = "'2.195"
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = "'0.3275"
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = "'0.1320"
$ws.Range("E26").Value = '25ProBitTokenPROB'
$ws.Range("D27").Value = "'0.0003157"
$ws.Range("D40").Value = "'0.03911"
$ws.Range("D41").Value = "'0.006760"
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("D42").Value = "'0.1072"
$ws.Range("D43").Value = "'0.002419"
$ws.Range("D44").Value = "'0.007753"
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
$ws.Range("D45").Value = "'0.00005632"
$ws.Range("D47").Value = "'0.3798"
$ws.Range("D48").Value = "'0.1684"
$ws.Range("D49").Value = "'0.00002099"
